$d = $word.ActiveDocument

$replacements = @(
    @("2024-08-14 Wednesday", "2024-08-15 Thursday"),
    @("35×36=", "85×14="),
    @("76×92=", "20×18="),
    @("23×86=", "67×81="),
    @("91×68=", "30×18="),
    @("80×67=", "97×41="),
    @("62×51=", "66×91="),
    @("19×97=", "72×26="),
    @("96×89=", "49×32="),
    @("56×70=", "93×81="),
    @("12×58=", "20×64="),
    @("99×34=", "65×42="),
    @("26×50=", "93×53="),
    @("56×99=", "47×55="),
    @("83×33=", "92×62="),
    @("43×69=", "89×63="),
    @("60×71=", "33×46="),
    @("49×29=", "17×67="),
    @("75×48=", "88×14="),
    @("61×26=", "19×25="),
    @("11×53=", "57×82="),
    @("48×31=", "97×91="),
    @("29×73=", "76×85="),
    @("54×46=", "49×24="),
    @("63×20=", "98×33="),
    @("13×25=", "44×93=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
